# This script applies the market-data refresh captured by the commit
# "chore: update Sheets via scheduled runner" to the Raiden_Profits workbook.
# It updates computed market-price / profit columns (H:N) on several of the
# per-job sheets, and clears the stale market-data columns (H:N) for the
# rows on the BSM sheet that no longer have pricing data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H99").Value = 783.6
$ws.Range("I99").Value = 783.6
$ws.Range("K99").Value = 2350.8
$ws.Range("M99").Value = -852.8000000000002

$ws.Range("H101").Value = 25316.111
$ws.Range("I101").Value = 20832.5
$ws.Range("K101").Value = 62497.5
$ws.Range("M101").Value = -60875.5

$ws.Range("H111").Value = 645
$ws.Range("I111").Value = 620.8333
$ws.Range("J111").Value = 790
$ws.Range("K111").Value = 1862.4999
$ws.Range("L111").Value = 2370
$ws.Range("M111").Value = 1204.5001
$ws.Range("N111").Value = -8504

$ws.Range("H115").Value = 1085.1111
$ws.Range("I115").Value = 1085.1111
$ws.Range("K115").Value = 3255.3333
$ws.Range("M115").Value = -1688.3333

$ws.Range("H118").Value = 322.16666
$ws.Range("I118").Value = 322.16666
$ws.Range("K118").Value = 966.4999799999999
$ws.Range("M118").Value = 690.5000200000001

$ws.Range("H127").Value = 1065.6
$ws.Range("I127").Value = 1135.1111
$ws.Range("J127").Value = 440
$ws.Range("K127").Value = 3405.3333
$ws.Range("L127").Value = 1320
$ws.Range("M127").Value = 1554.6667
$ws.Range("N127").Value = -11240

$ws.Range("H129").Value = 2777.1667
$ws.Range("I129").Value = 2953.8
$ws.Range("K129").Value = 8861.400000000001
$ws.Range("M129").Value = -3861.400000000001

$ws.Range("H132").Value = 32587.541
$ws.Range("I132").Value = 36985.953
$ws.Range("J132").Value = 1798.6666
$ws.Range("K132").Value = 110957.859
$ws.Range("L132").Value = 5395.9998
$ws.Range("M132").Value = -108427.859
$ws.Range("N132").Value = -10455.9998

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H61").Value = 4233.9165
$ws.Range("I61").Value = 3312.2222
$ws.Range("K61").Value = 3312.2222
$ws.Range("M61").Value = -3100.2222

$ws.Range("H97").Value = 536.1
$ws.Range("I97").Value = 530.35297
$ws.Range("J97").Value = 568.6667
$ws.Range("K97").Value = 530.35297
$ws.Range("L97").Value = 568.6667
$ws.Range("M97").Value = -34.35297000000003
$ws.Range("N97").Value = -1560.6667

$ws.Range("H132").Value = 3485.5715
$ws.Range("I132").Value = 3485.5715
$ws.Range("K132").Value = 10456.7145
$ws.Range("M132").Value = -7926.7145

$ws.Range("H136").Value = 4233.9165
$ws.Range("I136").Value = 3312.2222
$ws.Range("K136").Value = 9936.6666
$ws.Range("M136").Value = -7386.6666

# ---------------------------------------------------------------------
# BSM - clear stale market-data (columns H:N) for rows 117-141
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117:N141").ClearContents()

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 1720.1111
$ws.Range("J31").Value = 2310.5
$ws.Range("L31").Value = 2310.5
$ws.Range("N31").Value = -2900.5

$ws.Range("H34").Value = 1720.1111
$ws.Range("J34").Value = 2310.5
$ws.Range("L34").Value = 2310.5
$ws.Range("N34").Value = -2714.5

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 60000
$ws.Range("N96").Value = -64118

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 3891.3333
$ws.Range("I22").Value = 2400
$ws.Range("J22").Value = 4189.6
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 4189.6
$ws.Range("M22").Value = -2105
$ws.Range("N22").Value = -4779.6

$ws.Range("H27").Value = 3891.3333
$ws.Range("I27").Value = 2400
$ws.Range("J27").Value = 4189.6
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 4189.6
$ws.Range("M27").Value = -2293
$ws.Range("N27").Value = -4403.6

$ws.Range("H55").Value = 697.4
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 2000
$ws.Range("N55").Value = -2346

$ws.Range("H132").Value = 3481.3572
$ws.Range("J132").Value = 8500
$ws.Range("L132").Value = 25500
$ws.Range("N132").Value = -30560

# ---------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H116").Value = 109999.5
$ws.Range("J116").Value = 109999.5
$ws.Range("L116").Value = 109999.5
$ws.Range("N116").Value = -119177.5

$ws.Range("H126").Value = 2892.875
$ws.Range("I126").Value = 2892.875
$ws.Range("K126").Value = 8678.625
$ws.Range("M126").Value = -6208.625

$ws.Range("H132").Value = 6397.727
$ws.Range("I132").Value = 1239.4286
$ws.Range("K132").Value = 3718.2858
$ws.Range("M132").Value = -1188.2858

$ws.Range("H136").Value = 8124
$ws.Range("I136").Value = 7832
$ws.Range("K136").Value = 23496
$ws.Range("M136").Value = -20946
